$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B18").Formula = "=1000/0.91"
$ws.Range("B19").Select()
